$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '31.084.54'
$ws.Range('E2').Value = '  +2.90%  '
$ws.Range('D3').Value = '1.893.83'
$ws.Range('E3').Value = '  +2.95%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9965'
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.12'
$ws.Range('E5').Value = '  +2.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9959'
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4802'
$ws.Range('E7').Value = '  +2.65%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2865'
$ws.Range('E8').Value = '  +5.89%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06562'
$ws.Range('E9').Value = '  +4.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.87'
$ws.Range('E10').Value = '  +17.29%  '
$ws.Range('D11').Value = '1.930.65'
$ws.Range('E11').Value = '  +5.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '96.16'
$ws.Range('E12').Value = '  +14.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07532'
$ws.Range('E13').Value = '  +1.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.157'
$ws.Range('E14').Value = '  +4.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6568'
$ws.Range('E15').Value = '  +5.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '298.60'
$ws.Range('E16').Value = '  +31.48%  '
$ws.Range('D17').Value = '31.058.28'
$ws.Range('E17').Value = '  +3.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.18'
$ws.Range('E18').Value = '  +6.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9996'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007583'
$ws.Range('E20').Value = '  +3.93%  '
$ws.Range('B21').Value = 'BinanceUSD'
$ws.Range('C21').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9956'
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.213'
$ws.Range('E22').Value = '  +6.81%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.168'
$ws.Range('E23').Value = '  +5.52%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.342'
$ws.Range('E24').Value = '  +1.53%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '168.48'
$ws.Range('E25').Value = '  +2.64%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '19.69'
$ws.Range('E26').Value = '  +10.65%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.978'
$ws.Range('E27').Value = '  +5.36%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1055'
$ws.Range('E28').Value = '  +0.95%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.367'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.170'
$ws.Range('E30').Value = '  +2.36%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.002'
$ws.Range('E31').Value = '  +5.57%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05015'
$ws.Range('E32').Value = '  +4.19%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.190'
$ws.Range('E33').Value = '  +4.53%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7290'
$ws.Range('E34').Value = '  +2.77%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.697'
$ws.Range('E35').Value = '  +0.19%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.01947'
$ws.Range('E36').Value = '  +3.35%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.723'
$ws.Range('E37').Value = '  +2.84%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.076'
$ws.Range('E38').Value = '  +8.26%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9011'
$ws.Range('E39').Value = '  +1.14%  '
$ws.Range('B40').Value = 'Quant'
$ws.Range('C40').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '107.51'
$ws.Range('E40').Value = '  +3.18%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4247'
$ws.Range('E41').Value = '  +6.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9952'
$ws.Range('E42').Value = '  -0.63%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.595'
$ws.Range('E43').Value = '  +1.12%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.454'
$ws.Range('E44').Value = '  +6.09%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.94'
$ws.Range('E45').Value = '  +10.34%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1239'
$ws.Range('E46').Value = '  +3.80%  '
$ws.Range('B47').Value = 'Elrond'
$ws.Range('C47').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '34.76'
$ws.Range('E47').Value = '  +6.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.907'
$ws.Range('E48').Value = '  +4.76%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.406'
$ws.Range('E49').Value = '  +3.80%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05584'
$ws.Range('E50').Value = '  +1.36%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3825'
$ws.Range('E51').Value = '  +5.45%  '
